# Deploy the implementation guide.
#
# - rename the "Include from Ferlab.bio CodeS" sheet to "Include #0"
# - refresh the Metadata sheet's Date value
# - refresh the Metadata sheet's Contact value
# - add a new "Jurisdiction" row (empty value) right after "Contact",
#   pushing Description / Purpose / Copyright / Immutable down by one row

$wb = $excel.ActiveWorkbook

# --- Rename the second (Include) worksheet -------------------------------
$includeSheet = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$includeSheet.Name = "Include #0"

# --- Metadata worksheet updates -------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Date (row 8, column B)
$meta.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Contact (row 10, column B)
$meta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row for "Jurisdiction" right below "Contact" (row 11),
# copying the existing row formatting so the new row matches the table
# style, then overwrite its values.
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy($meta.Range("A11:B11"))
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""
